# Update the "predicted_up" column (column C) values for specific rows
# to reflect the refreshed model predictions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 0
    5  = 1
    16 = 1
    27 = 0
    33 = 0
    45 = 1
    46 = 1
    50 = 1
    51 = 1
    58 = 1
    59 = 1
    65 = 0
    73 = 1
    74 = 1
    86 = 0
    90 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}
